$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell-level value / missing-data mask changes (applied using original row numbers) ---
$ws.Range("D2").Value = -13.5    # RM 2
$ws.Range("F3").ClearContents()  # RM 8
$ws.Range("F4").Value = 17.97    # RM 9
$ws.Range("F5").ClearContents()  # RM 14
$ws.Range("D6").ClearContents()  # RM 21
$ws.Range("F8").ClearContents()  # RM 38
$ws.Range("D12").Value = -14.1   # RM 81
$ws.Range("D14").ClearContents() # RM 90
$ws.Range("D20").Value = -14.0   # RM 134
$ws.Range("D21").Value = -14.3   # RM 135
$ws.Range("D22").ClearContents() # RM 138
$ws.Range("D23").ClearContents() # RM 140
$ws.Range("F23").Value = 16.48   # RM 140
$ws.Range("F29").ClearContents() # SC 101
$ws.Range("F31").Value = 18.06   # SC 119
$ws.Range("C32").Value = 11.4    # SC 120
$ws.Range("D33").Value = -13.7   # SC 132
$ws.Range("C34").ClearContents() # SC 193
$ws.Range("D35").Value = -14.1   # SC 232

# --- Remove the two rows that are no longer present (delete higher index first) ---
$ws.Rows("28").Delete()  # SC 92
$ws.Rows("26").Delete()  # RM 232
